$d = $word.ActiveDocument

# Locate the "2PE" Heading 2 paragraph, then remove the very next paragraph,
# which holds the standalone italic "2 Pierre" run (whole paragraph, incl.
# its paragraph mark, is being dropped).
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd() -eq "2PE" -and $p.Style.NameLocal -eq "Heading 2") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the '2PE' Heading 2 paragraph"
}

$italicParaIndex = $targetIndex + 1
$italicPara = $d.Paragraphs.Item($italicParaIndex)

if ($italicPara.Range.Text.TrimEnd() -ne "2 Pierre") {
    throw "Unexpected paragraph after '2PE' heading: '$($italicPara.Range.Text)'"
}

$italicPara.Range.Delete()

Write-Output "Removed standalone '2 Pierre' paragraph after the 2PE heading."
